# v1.9: Comisiones solo cat. Comisiones; Total flujo por mes; Errores sin Editado, scroll, Descripcion izq.
# Appends one new log entry to the "Log" sheet (row 40) and one new version
# entry to the "Versiones" sheet (row 11).

$wb = $excel.ActiveWorkbook

# --- Sheet "Log": add row 40 -------------------------------------------
$logSheet = $wb.Worksheets.Item("Log")

$logSheet.Cells.Item(40, 1).Value = "27/02/2025"
$logSheet.Cells.Item(40, 2).Value = "20:10"
$logSheet.Cells.Item(40, 3).Value = "Comisiones, Total flujo, Errores"
$logSheet.Cells.Item(40, 4).Value = "Comisiones solo categoría Comisiones (no Sueldos). Fila Total en flujo por mes con sumas y ratios. Tabla Errores: quitar columna Editado, scroll horizontal, Descripción alineada a la izquierda."
$logSheet.Cells.Item(40, 5).Value = "Diagnostico"

# --- Sheet "Versiones": add row 11 --------------------------------------
$verSheet = $wb.Worksheets.Item("Versiones")

# "1.9" looks numeric, so a plain .Value assignment would be stored as the
# number 1.9 instead of the text "1.9" (every other cell in this workbook,
# including the existing version numbers, is stored as text). Entering it
# as a formula that yields the text "1.9" and then converting that formula
# to its resulting value (copy / paste-special-values) keeps the cell a
# plain text cell without touching the sheet's cell-formatting table.
$verSheet.Range("A11").Formula = '="1.9"'
$verSheet.Range("A11").Copy()
$verSheet.Range("A11").PasteSpecial(-4163)

$verSheet.Cells.Item(11, 2).Value = "27/02/2025"
$verSheet.Cells.Item(11, 3).Value = "Comisiones solo categoría Comisiones; fila Total en flujo por mes; Errores: sin columna Editado, scroll horizontal, Descripción a la izquierda"
